$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts old B->C, old C->D)
$ws.Columns.Item(2).Insert()

# Set width of new column B to match column A
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

# Header row
$ws.Range("B1").Value = "StatQuery"

# Update A2 (diagnosis query) text to the new disease term
$ws.Range("A2").Value = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN ['Bone sarcomas :: Osteosarcoma (appendicular)'] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'') AS ``Case ID`` , coalesce(s.clinical_study_designation,'') AS ``Study Code`` , coalesce(s.clinical_study_type,'') AS  ``Study Type``, coalesce(demo.breed,'') AS Breed , coalesce(diag.disease_term,'') AS Diagnosis , coalesce(diag.stage_of_disease,'') AS ``Stage of Disease`` ,  coalesce(demo.patient_age_at_enrollment,'') AS Age , coalesce(demo.sex,'') AS Sex , coalesce(demo.neutered_indicator,'') AS  ``Neutered Status``"

# Set B2 (new stat query) with same style as A2 (wrap text)
$ws.Range("B2").Value = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN['Bone sarcomas :: Osteosarcoma (appendicular)']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"
$ws.Range("B2").WrapText = $true

# Update sheet view/selection
$ws.Application.ActiveWindow.SelectedSheets.Item(1).Select()
$ws.Range("B1:B1048576").Select()
